$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 260.1111
$ws.Range("I9").Value = 280.125
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 280.125
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -111.125
$ws.Range("N9").Value = -438

$ws.Range("H39").Value = 1443171.6
$ws.Range("I39").Value = 1763787.5
$ws.Range("J39").Value = 400
$ws.Range("K39").Value = 5291362.5
$ws.Range("L39").Value = 1200
$ws.Range("M39").Value = -5291066.5
$ws.Range("N39").Value = -1792

$ws.Range("H62").Value = 2270.7144
$ws.Range("I62").Value = 2270.7144
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2270.7144
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1646.7144

$ws.Range("H65").Value = 2270.7144
$ws.Range("I65").Value = 2270.7144
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11353.572
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -8233.572

$ws.Range("H106").Value = 2525.5
$ws.Range("I106").Value = 2672
$ws.Range("J106").Value = 1500
$ws.Range("K106").Value = 2672
$ws.Range("L106").Value = 1500
$ws.Range("M106").Value = -2041
$ws.Range("N106").Value = -2762

$ws.Range("H112").Value = 1057.2609
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 1059.8636
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 3179.5908
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -5395.5908

$ws.Range("H125").Value = 2423.2
$ws.Range("I125").Value = 4310.6665
$ws.Range("J125").Value = 1614.2858
$ws.Range("K125").Value = 38795.9985
$ws.Range("L125").Value = 14528.5722
$ws.Range("M125").Value = -36335.9985
$ws.Range("N125").Value = -19448.5722

$ws.Range("H132").Value = 5957844.5
$ws.Range("I132").Value = 6584732
$ws.Range("J132").Value = 2415.25
$ws.Range("K132").Value = 19754196
$ws.Range("L132").Value = 7245.75
$ws.Range("M132").Value = -19751666
$ws.Range("N132").Value = -12305.75

$ws.Range("H137").Value = 1404.9474
$ws.Range("I137").Value = 1171.9166
$ws.Range("J137").Value = 1574.4242
$ws.Range("K137").Value = 3515.7498
$ws.Range("L137").Value = 4723.2726
$ws.Range("M137").Value = -965.7498000000001
$ws.Range("N137").Value = -9823.2726

$ws.Range("H138").Value = 5878.079
$ws.Range("I138").Value = 4249.625
$ws.Range("J138").Value = 6312.3335
$ws.Range("K138").Value = 12748.875
$ws.Range("L138").Value = 18937.0005
$ws.Range("M138").Value = -7608.875
$ws.Range("N138").Value = -29217.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21762.012
$ws.Range("I32").Value = 3952.5
$ws.Range("J32").Value = 87063.55499999999
$ws.Range("K32").Value = 3952.5
$ws.Range("L32").Value = 87063.55499999999
$ws.Range("M32").Value = -3665.5
$ws.Range("N32").Value = -87637.55499999999

$ws.Range("H74").Value = 3019.7585
$ws.Range("I74").Value = 2412.1667
$ws.Range("J74").Value = 3448.647
$ws.Range("K74").Value = 2412.1667
$ws.Range("L74").Value = 3448.647
$ws.Range("M74").Value = -1538.1667
$ws.Range("N74").Value = -5196.647

$ws.Range("H77").Value = 3019.7585
$ws.Range("I77").Value = 2412.1667
$ws.Range("J77").Value = 3448.647
$ws.Range("K77").Value = 12060.8335
$ws.Range("L77").Value = 17243.235
$ws.Range("M77").Value = -7692.833500000001
$ws.Range("N77").Value = -25979.235

$ws.Range("H88").Value = 2265.889
$ws.Range("I88").Value = 1900
$ws.Range("J88").Value = 2448.8333
$ws.Range("K88").Value = 1900
$ws.Range("L88").Value = 2448.8333
$ws.Range("M88").Value = -1494
$ws.Range("N88").Value = -3260.8333

$ws.Range("H91").Value = 2265.889
$ws.Range("I91").Value = 1900
$ws.Range("J91").Value = 2448.8333
$ws.Range("K91").Value = 1900
$ws.Range("L91").Value = 2448.8333
$ws.Range("M91").Value = -496
$ws.Range("N91").Value = -5256.8333

$ws.Range("H122").Value = 2402.8462
$ws.Range("I122").Value = 2488.2307
$ws.Range("J122").Value = 2232.077
$ws.Range("K122").Value = 7464.6921
$ws.Range("L122").Value = 6696.231000000001
$ws.Range("M122").Value = -5014.6921
$ws.Range("N122").Value = -11596.231

$ws.Range("H132").Value = 3148.125
$ws.Range("I132").Value = 3961.2415
$ws.Range("J132").Value = 1907.0526
$ws.Range("K132").Value = 11883.7245
$ws.Range("L132").Value = 5721.1578
$ws.Range("M132").Value = -9353.7245
$ws.Range("N132").Value = -10781.1578

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 37153.105
$ws.Range("I20").Value = 60506.06
$ws.Range("J20").Value = 1062.1818
$ws.Range("K20").Value = 60506.06
$ws.Range("L20").Value = 1062.1818
$ws.Range("M20").Value = -60259.06
$ws.Range("N20").Value = -1556.1818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14169.21
$ws.Range("I31").Value = 28882.584
$ws.Range("J31").Value = 2398.5112
$ws.Range("K31").Value = 28882.584
$ws.Range("L31").Value = 2398.5112
$ws.Range("M31").Value = -28587.584
$ws.Range("N31").Value = -2988.5112

$ws.Range("H34").Value = 14169.21
$ws.Range("I34").Value = 28882.584
$ws.Range("J34").Value = 2398.5112
$ws.Range("K34").Value = 28882.584
$ws.Range("L34").Value = 2398.5112
$ws.Range("M34").Value = -28680.584
$ws.Range("N34").Value = -2802.5112

$ws.Range("H132").Value = 3025.8147
$ws.Range("I132").Value = 2966.6191
$ws.Range("J132").Value = 3233
$ws.Range("K132").Value = 8899.8573
$ws.Range("L132").Value = 9699
$ws.Range("M132").Value = -6369.8573
$ws.Range("N132").Value = -14759

$ws.Range("H134").Value = 3316.2
$ws.Range("I134").Value = 2366.8
$ws.Range("J134").Value = 4265.6
$ws.Range("K134").Value = 7100.400000000001
$ws.Range("L134").Value = 12796.8
$ws.Range("M134").Value = -4565.400000000001
$ws.Range("N134").Value = -17866.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 521.6
$ws.Range("I14").Value = 521.6
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1564.8
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1391.8

$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2712

$ws.Range("H44").Value = 900.6
$ws.Range("I44").Value = 499.5
$ws.Range("J44").Value = 1168
$ws.Range("K44").Value = 1498.5
$ws.Range("L44").Value = 3504
$ws.Range("M44").Value = -1100.5
$ws.Range("N44").Value = -4300

$ws.Range("H113").Value = 548.3182
$ws.Range("I113").Value = 535.73334
$ws.Range("J113").Value = 554.8276
$ws.Range("K113").Value = 1607.20002
$ws.Range("L113").Value = 1664.4828
$ws.Range("M113").Value = 562.79998
$ws.Range("N113").Value = -6004.4828

$ws.Range("H137").Value = 2600
$ws.Range("I137").Value = 2137.7778
$ws.Range("J137").Value = 2920
$ws.Range("K137").Value = 6413.3334
$ws.Range("L137").Value = 8760
$ws.Range("M137").Value = -1313.3334
$ws.Range("N137").Value = -18960

$ws.Range("H138").Value = 1891.3334
$ws.Range("I138").Value = 1317
$ws.Range("J138").Value = 3040
$ws.Range("K138").Value = 3951
$ws.Range("L138").Value = 9120
$ws.Range("M138").Value = 1189
$ws.Range("N138").Value = -19400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 10472.272
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 10472.272
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 10472.272
$ws.Range("N46").Value = -10784.272

$ws.Range("H122").Value = 5364.591
$ws.Range("I122").Value = 5683.75
$ws.Range("J122").Value = 4981.6
$ws.Range("K122").Value = 17051.25
$ws.Range("L122").Value = 14944.8
$ws.Range("M122").Value = -14601.25
$ws.Range("N122").Value = -19844.8

$ws.Range("H132").Value = 2222.238
$ws.Range("I132").Value = 1580.0769
$ws.Range("J132").Value = 3265.75
$ws.Range("K132").Value = 4740.2307
$ws.Range("L132").Value = 9797.25
$ws.Range("M132").Value = -2210.2307
$ws.Range("N132").Value = -14857.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5333.278
$ws.Range("I132").Value = 6650.2
$ws.Range("J132").Value = 3687.125
$ws.Range("K132").Value = 19950.6
$ws.Range("L132").Value = 11061.375
$ws.Range("M132").Value = -17420.6
$ws.Range("N132").Value = -16121.375

$ws.Range("H136").Value = 2088.5
$ws.Range("I136").Value = 1604.591
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 4813.772999999999
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -2263.772999999999
$ws.Range("N136").Value = -19350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 6786.4614
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 6786.4614
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 6786.4614
$ws.Range("N49").Value = -7246.4614

$ws.Range("H107").Value = 40999.668
$ws.Range("I107").Value = 10630.85
$ws.Range("J107").Value = 101737.3
$ws.Range("K107").Value = 31892.55
$ws.Range("L107").Value = 305211.9
$ws.Range("M107").Value = -29972.55
$ws.Range("N107").Value = -309051.9

$ws.Range("H113").Value = 275.92307
$ws.Range("I113").Value = 199.33333
$ws.Range("J113").Value = 448.25
$ws.Range("K113").Value = 597.99999
$ws.Range("L113").Value = 1344.75
$ws.Range("M113").Value = 1572.00001
$ws.Range("N113").Value = -5684.75

$ws.Range("H122").Value = 1353.5385
$ws.Range("I122").Value = 1515.7142
$ws.Range("J122").Value = 1164.3334
$ws.Range("K122").Value = 4547.142599999999
$ws.Range("L122").Value = 3493.0002
$ws.Range("M122").Value = -2097.142599999999
$ws.Range("N122").Value = -8393.0002

$ws.Range("H126").Value = 2238.3333
$ws.Range("I126").Value = 3143.3333
$ws.Range("J126").Value = 1333.3334
$ws.Range("K126").Value = 9429.999899999999
$ws.Range("L126").Value = 4000.0002
$ws.Range("M126").Value = -6959.999899999999
$ws.Range("N126").Value = -8940.0002

$ws.Range("H136").Value = 1307.5714
$ws.Range("I136").Value = 777.26666
$ws.Range("J136").Value = 2633.3333
$ws.Range("K136").Value = 2331.79998
$ws.Range("L136").Value = 7899.999899999999
$ws.Range("M136").Value = 218.2000200000002
$ws.Range("N136").Value = -12999.9999
